$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 11422
$ws.Range("F4").Value = 1304
$ws.Range("F6").Value = 1235
$ws.Range("G6").Value = 60
$ws.Range("F9").Value = 963
$ws.Range("F11").Value = 2327
$ws.Range("F13").Value = 1143
$ws.Range("F14").Value = 877
$ws.Range("F15").Value = 583
$ws.Range("F16").Value = 866
$ws.Range("F17").Value = 1038
$ws.Range("F21").Value = 725
$ws.Range("F22").Value = 156
$ws.Range("F23").Value = 425
$ws.Range("F24").Value = 1073
$ws.Range("F26").Value = 489
$ws.Range("F27").Value = 550
$ws.Range("F30").Value = 277
$ws.Range("F31").Value = 644
$ws.Range("F32").Value = 2728
$ws.Range("F33").Value = 442
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 292
$ws.Range("F37").Value = 93
$ws.Range("F38").Value = 1530
$ws.Range("F39").Value = 433
$ws.Range("F41").Value = 68
$ws.Range("F42").Value = 117
$ws.Range("F47").Value = 68
$ws.Range("F48").Value = 15
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 18
$ws.Range("F7").Value = 94
$ws.Range("F10").Value = 163
$ws.Range("F11").Value = 4408
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 125
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2233
$ws.Range("F3").Value = 687
$ws.Range("F4").Value = 651
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2233
$ws.Range("F4").Value = 11422
$ws.Range("F5").Value = 651
$ws.Range("F6").Value = 1235
$ws.Range("G6").Value = 60
$ws.Range("F7").Value = 18
$ws.Range("F10").Value = 963
$ws.Range("F11").Value = 2327
$ws.Range("F13").Value = 1143
$ws.Range("F14").Value = 877
$ws.Range("F15").Value = 583
$ws.Range("F16").Value = 866
$ws.Range("F17").Value = 1038
$ws.Range("F24").Value = 725
$ws.Range("F25").Value = 156
$ws.Range("F26").Value = 425
$ws.Range("F27").Value = 1073
$ws.Range("F28").Value = 94
$ws.Range("F30").Value = 489
$ws.Range("F31").Value = 550
$ws.Range("F34").Value = 2728
$ws.Range("F35").Value = 163
$ws.Range("F36").Value = 442
$ws.Range("F37").Value = 93
$ws.Range("F38").Value = 1530
$ws.Range("F39").Value = 433
$ws.Range("F42").Value = 15
$ws.Range("F43").Value = 118
$ws.Range("F47").Value = 68
